# Applies the "Updated cryptos list" price/volume refresh described by the diff.
# D-column values are text-formatted numbers (e.g. "27.863.77"), so they are
# entered with a leading apostrophe to force text entry, matching the original
# inline-string cell contents instead of letting Excel reinterpret them as numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "D2"; Value = "'27.863.77" },
    @{ Cell = "E2"; Value = "  -0.29%  " },
    @{ Cell = "D3"; Value = "'1.626.61" },
    @{ Cell = "E3"; Value = "  -0.95%  " },
    @{ Cell = "E4"; Value = "  -0.33%  " },
    @{ Cell = "D5"; Value = "'211.18" },
    @{ Cell = "E5"; Value = "  -0.79%  " },
    @{ Cell = "D6"; Value = "'0.521" },
    @{ Cell = "E6"; Value = "  -0.37%  " },
    @{ Cell = "E7"; Value = "  -0.33%  " },
    @{ Cell = "D8"; Value = "'23.32" },
    @{ Cell = "E8"; Value = "  -0.63%  " },
    @{ Cell = "E9"; Value = "  -0.68%  " },
    @{ Cell = "D10"; Value = "'0.0611" },
    @{ Cell = "E10"; Value = "  -0.43%  " },
    @{ Cell = "D11"; Value = "'0.0878" },
    @{ Cell = "E11"; Value = "  -0.57%  " },
    @{ Cell = "D12"; Value = "'1.857.49" },
    @{ Cell = "E12"; Value = "  -0.93%  " },
    @{ Cell = "D13"; Value = "'1.625.91" },
    @{ Cell = "E13"; Value = "  -1.02%  " },
    @{ Cell = "E14"; Value = "  -1.63%  " },
    @{ Cell = "D15"; Value = "'0.562" },
    @{ Cell = "E15"; Value = "  -1.85%  " },
    @{ Cell = "D16"; Value = "'65.26" },
    @{ Cell = "E16"; Value = "  -0.48%  " },
    @{ Cell = "D17"; Value = "'27.858.51" },
    @{ Cell = "E17"; Value = "  -0.29%  " },
    @{ Cell = "D18"; Value = "'229.63" },
    @{ Cell = "E18"; Value = "  -1.36%  " },
    @{ Cell = "D20"; Value = "'0.0₃0720" },
    @{ Cell = "E20"; Value = "  -0.38%  " },
    @{ Cell = "E21"; Value = "  -0.33%  " },
    @{ Cell = "E22"; Value = "  -1.27%  " },
    @{ Cell = "D23"; Value = "'10.07" },
    @{ Cell = "E23"; Value = "  -3.71%  " },
    @{ Cell = "D24"; Value = "'2.05" },
    @{ Cell = "E24"; Value = "  -2.08%  " },
    @{ Cell = "D25"; Value = "'154.39" },
    @{ Cell = "E25"; Value = "  +0.97%  " },
    @{ Cell = "D26"; Value = "'6.89" },
    @{ Cell = "E26"; Value = "  -0.13%  " },
    @{ Cell = "E27"; Value = "  -0.18%  " },
    @{ Cell = "D28"; Value = "'15.52" },
    @{ Cell = "E28"; Value = "  -1.15%  " },
    @{ Cell = "D29"; Value = "'0.998" },
    @{ Cell = "E29"; Value = "  -0.23%  " },
    @{ Cell = "D30"; Value = "'1.17" },
    @{ Cell = "E30"; Value = "  -1.26%  " },
    @{ Cell = "E31"; Value = "  -0.85%  " },
    @{ Cell = "E32"; Value = "  +1.56%  " },
    @{ Cell = "E33"; Value = "  -0.41%  " },
    @{ Cell = "D34"; Value = "'1.397.05" },
    @{ Cell = "E34"; Value = "  -0.78%  " },
    @{ Cell = "E35"; Value = "  +0.31%  " },
    @{ Cell = "E36"; Value = "  +10.58%  " },
    @{ Cell = "E37"; Value = "  -1.02%  " },
    @{ Cell = "E38"; Value = "  +0.02%  " },
    @{ Cell = "D39"; Value = "'0.554" },
    @{ Cell = "E39"; Value = "  -1.46%  " },
    @{ Cell = "D40"; Value = "'0.851" },
    @{ Cell = "E40"; Value = "  -3.51%  " },
    @{ Cell = "E41"; Value = "  -0.30%  " },
    @{ Cell = "D42"; Value = "'1.01" },
    @{ Cell = "E42"; Value = "  -1.89%  " },
    @{ Cell = "E43"; Value = "  -0.20%  " },
    @{ Cell = "D44"; Value = "'65.66" },
    @{ Cell = "E44"; Value = "  -2.46%  " },
    @{ Cell = "D45"; Value = "'5.42" },
    @{ Cell = "E45"; Value = "  -1.76%  " },
    @{ Cell = "D46"; Value = "'1.768.99" },
    @{ Cell = "E46"; Value = "  -0.78%  " },
    @{ Cell = "E47"; Value = "  -3.02%  " },
    @{ Cell = "E48"; Value = "  +0.28%  " },
    @{ Cell = "E49"; Value = "  +1.36%  " },
    @{ Cell = "D50"; Value = "'0.0₆0100" },
    @{ Cell = "E50"; Value = "  +10.16%  " },
    @{ Cell = "D51"; Value = "'0.0502" },
    @{ Cell = "E51"; Value = "  -0.77%  " }
)

foreach ($u in $updates) {
    $ws.Range($u.Cell).Value = $u.Value
}
